# Added more info slide at end
$p = $ppt.ActivePresentation

# Append a new slide, using the same "Title and Content" layout (index 2,
# ppLayoutText) used by the preceding slide.
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "For More Information:"

# Body / content placeholder - build it up as several runs the way PowerPoint
# does when a URL is typed and auto-corrected/spell-split on its word
# boundaries.
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "https://"
$tr.InsertAfter("github.com") | Out-Null
$tr.InsertAfter("/") | Out-Null
$tr.InsertAfter("HumanDynamics") | Out-Null
$tr.InsertAfter("/") | Out-Null
$tr.InsertAfter("CodeTheDeal") | Out-Null
